# Update TPM-derived NATMI ligand-receptor stats (Psen1-Notch3) with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value2 = 9.409481333333334
$ws.Cells.Item(2,8).Value2 = 28.228444
$ws.Cells.Item(2,9).Value2 = 0.2433300530093958
$ws.Cells.Item(2,10).Value2 = 0.2433300530093958
$ws.Cells.Item(2,13).Value2 = 7.413580666666667
$ws.Cells.Item(2,14).Value2 = 22.240742
$ws.Cells.Item(2,15).Value2 = 0.05108888817597561
$ws.Cells.Item(2,16).Value2 = 0.05108888817597561
$ws.Cells.Item(2,17).Value2 = 69.7579488961609
$ws.Cells.Item(2,18).Value2 = 627.8215400654481
$ws.Cells.Item(2,19).Value2 = 0.01243146186805124
$ws.Cells.Item(2,20).Value2 = 0.01243146186805124
$ws.Cells.Item(3,7).Value2 = 9.409481333333334
$ws.Cells.Item(3,8).Value2 = 28.228444
$ws.Cells.Item(3,9).Value2 = 0.2433300530093958
$ws.Cells.Item(3,10).Value2 = 0.2433300530093958
$ws.Cells.Item(3,15).Value2 = 0.0112127179963522
$ws.Cells.Item(3,16).Value2 = 0.0112127179963522
$ws.Cells.Item(3,17).Value2 = 15.31010434759111
$ws.Cells.Item(3,18).Value2 = 137.79093912832
$ws.Cells.Item(3,19).Value2 = 0.002728391264431787
$ws.Cells.Item(3,20).Value2 = 0.002728391264431788
$ws.Cells.Item(4,7).Value2 = 9.409481333333334
$ws.Cells.Item(4,8).Value2 = 28.228444
$ws.Cells.Item(4,9).Value2 = 0.2433300530093958
$ws.Cells.Item(4,10).Value2 = 0.2433300530093958
$ws.Cells.Item(4,15).Value2 = 0.9376983938276722
$ws.Cells.Item(4,16).Value2 = 0.9376983938276722
$ws.Cells.Item(4,17).Value2 = 1280.355062950903
$ws.Cells.Item(4,18).Value2 = 11523.19556655813
$ws.Cells.Item(4,19).Value2 = 0.2281701998769128
$ws.Cells.Item(4,20).Value2 = 0.2281701998769128
$ws.Cells.Item(5,9).Value2 = 0.5069354697952918
$ws.Cells.Item(5,10).Value2 = 0.5069354697952919
$ws.Cells.Item(5,13).Value2 = 7.413580666666667
$ws.Cells.Item(5,14).Value2 = 22.240742
$ws.Cells.Item(5,15).Value2 = 0.05108888817597561
$ws.Cells.Item(5,16).Value2 = 0.05108888817597561
$ws.Cells.Item(5,17).Value2 = 145.3284465206022
$ws.Cells.Item(5,18).Value2 = 1307.95601868542
$ws.Cells.Item(5,19).Value2 = 0.02589876952880733
$ws.Cells.Item(5,20).Value2 = 0.02589876952880733
$ws.Cells.Item(6,9).Value2 = 0.5069354697952918
$ws.Cells.Item(6,10).Value2 = 0.5069354697952919
$ws.Cells.Item(6,15).Value2 = 0.0112127179963522
$ws.Cells.Item(6,16).Value2 = 0.0112127179963522
$ws.Cells.Item(6,19).Value2 = 0.005684124465162926
$ws.Cells.Item(6,20).Value2 = 0.005684124465162927
$ws.Cells.Item(7,9).Value2 = 0.5069354697952918
$ws.Cells.Item(7,10).Value2 = 0.5069354697952919
$ws.Cells.Item(7,15).Value2 = 0.9376983938276722
$ws.Cells.Item(7,16).Value2 = 0.9376983938276722
$ws.Cells.Item(7,19).Value2 = 0.4753525758013216
$ws.Cells.Item(7,20).Value2 = 0.4753525758013217
$ws.Cells.Item(8,7).Value2 = 9.657138
$ws.Cells.Item(8,9).Value2 = 0.2497344771953123
$ws.Cells.Item(8,10).Value2 = 0.2497344771953124
$ws.Cells.Item(8,13).Value2 = 7.413580666666667
$ws.Cells.Item(8,14).Value2 = 22.240742
$ws.Cells.Item(8,15).Value2 = 0.05108888817597561
$ws.Cells.Item(8,16).Value2 = 0.05108888817597561
$ws.Cells.Item(8,17).Value2 = 71.593971572132
$ws.Cells.Item(8,18).Value2 = 644.345744149188
$ws.Cells.Item(8,19).Value2 = 0.01275865677911704
$ws.Cells.Item(8,20).Value2 = 0.01275865677911704
$ws.Cells.Item(9,7).Value2 = 9.657138
$ws.Cells.Item(9,9).Value2 = 0.2497344771953123
$ws.Cells.Item(9,10).Value2 = 0.2497344771953124
$ws.Cells.Item(9,15).Value2 = 0.0112127179963522
$ws.Cells.Item(9,16).Value2 = 0.0112127179963522
$ws.Cells.Item(9,19).Value2 = 0.002800202266757487
$ws.Cells.Item(9,20).Value2 = 0.002800202266757487
$ws.Cells.Item(10,7).Value2 = 9.657138
$ws.Cells.Item(10,9).Value2 = 0.2497344771953123
$ws.Cells.Item(10,10).Value2 = 0.2497344771953124
$ws.Cells.Item(10,15).Value2 = 0.9376983938276722
$ws.Cells.Item(10,16).Value2 = 0.9376983938276722
$ws.Cells.Item(10,19).Value2 = 0.2341756181494378
$ws.Cells.Item(10,20).Value2 = 0.2341756181494378
